# Apply the benchmark table updates described in the commit:
# "Fixed README.md stats and docx preparation for all DaCapo - JDK 21 - Z GC tests"
#
# The document is a single-column, single-cell-per-row table. We address
# cells positionally (by row index) rather than by Find/Replace text, since
# several of the original values are not unique within the document (e.g.
# multiple cells contain "0.00004"), which would make a text-based search
# ambiguous or error-prone.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell corrections (rows 1-12)
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "102"
$t.Cell(5, 1).Range.Text  = "0.00002"
$t.Cell(6, 1).Range.Text  = "0.01046"
$t.Cell(7, 1).Range.Text  = "0.00153"
$t.Cell(8, 1).Range.Text  = "0.00000"
$t.Cell(9, 1).Range.Text  = "0.01046"
$t.Cell(10, 1).Range.Text = "0.01046"
$t.Cell(11, 1).Range.Text = "0.01046"
$t.Cell(12, 1).Range.Text = "0.01407"

# Rows 44-46 previously held a whole tab-separated raw data row crammed into
# a single run; collapse each back down to just its leading/summary value.
$t.Cell(44, 1).Range.Text = "99.98"
$t.Cell(45, 1).Range.Text = "0.01"
$t.Cell(46, 1).Range.Text = "66"
